$wb = $excel.ActiveWorkbook

# OFF sheet ("OFF") - Row 3 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 346
$wsOff.Range("C3").Value = 221
$wsOff.Range("D3").Value = 122
$wsOff.Range("E3").Value = 58

# DEF sheet ("DEF") - Row 3 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 453
$wsDef.Range("C3").Value = 297
$wsDef.Range("D3").Value = 120
$wsDef.Range("E3").Value = 45
$wsDef.Range("G3").Value = 7
